$wb = $excel.ActiveWorkbook

$wsPed = $wb.Worksheets.Item('Pediatric Vaccine ')
$wsAdult = $wb.Worksheets.Item('Adult Vaccine ')
$wsPedFlu = $wb.Worksheets.Item('Pediatric Influenza Vaccine ')
$wsAdultFlu = $wb.Worksheets.Item('Adult Influenza Vaccine ')

$wsPed.Range('A2').Value = 'DTaP '
$wsPed.Range('A3').Value = 'DTaP '
$wsPed.Range('A4').Value = 'DTaP '
$wsPed.Range('A5').Value = 'DTaP-IPV '
$wsPed.Range('A6').Value = 'DTaP-IPV '
$wsPed.Range('A7').Value = 'DTaP-Hep B-IPV '
$wsPed.Range('A8').Value = 'DTaP-IP-HI '
$wsPed.Range('A9').Value = 'e-IPV '
$wsPed.Range('A10').Value = 'Hepatitis A Pediatric '
$wsPed.Range('A11').Value = 'Hepatitis A Pediatric '
$wsPed.Range('A12').Value = 'Hepatitis A Pediatric '
$wsPed.Range('A13').Value = 'Hepatitis A Pediatric '
$wsPed.Range('A14').Value = 'Hepatitis A-Hepatitis B 18 only '
$wsPed.Range('A15').Value = 'Hepatitis A-Hepatitis B 18 only '
$wsPed.Range('A16').Value = 'Hepatitis B  Pediatric/Adolescent'
$wsPed.Range('A17').Value = 'Hepatitis B  Pediatric/Adolescent'
$wsPed.Range('A18').Value = 'Hepatitis B  Pediatric/Adolescent'
$wsPed.Range('A19').Value = 'Hepatitis B  Pediatric/Adolescent'
$wsPed.Range('B18').Value = 'Recombivax HB'
$wsPed.Range('B19').Value = 'Recombivax HB'
$wsPed.Range('A20').Value = 'Hib '
$wsPed.Range('A21').Value = 'Hib '
$wsPed.Range('A22').Value = 'HIBMENCY '
$wsPed.Range('A23').Value = 'HPV - Quadrivalent Human Papillomavirus Types 6, 11, 16 and 18 Recombinant '
$wsPed.Range('A24').Value = 'HPV - Human Papillomavirus 9-valent '
$wsPed.Range('A25').Value = 'HPV -Bivalent Human Papillomavirus Types 16 and 18 '
$wsPed.Range('A26').Value = 'MENB - Meningococcal Group B '
$wsPed.Range('A27').Value = 'MENB - Meningococcal Group B '
$wsPed.Range('A28').Value = 'MENB - Meningococcal Group B '
$wsAdult.Range('A18').Value = 'MENB - Meningococcal Group B '
$wsAdult.Range('A19').Value = 'MENB - Meningococcal Group B '
$wsAdult.Range('A20').Value = 'MENB - Meningococcal Group B '
$wsPed.Range('A29').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$wsPed.Range('A30').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$wsPed.Range('A31').Value = 'Measles, Mumps and Rubella (MMR) '
$wsPed.Range('A32').Value = 'MMR/Varicella '
$wsPed.Range('A33').Value = 'Pneumococcal 13-valent  (Pediatric)'
$wsPed.Range('A35').Value = 'Rotavirus, Live, Oral, Pentavalent '
$wsPed.Range('A36').Value = 'Rotavirus, Live, Oral, Pentavalent '
$wsPed.Range('A37').Value = 'Rotavirus, Live, Oral, Oral '
$wsPed.Range('A38').Value = 'Tetanus  Diphtheria Toxoids '
$wsPed.Range('A39').Value = 'Tetanus  Diphtheria Toxoids '
$wsPed.Range('A40').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsPed.Range('A41').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsPed.Range('A42').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsPed.Range('A43').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsAdult.Range('A26').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsAdult.Range('A27').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsAdult.Range('A28').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsAdult.Range('A29').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis '
$wsPed.Range('A44').Value = 'Varicella '
$wsAdult.Range('A2').Value = 'Hepatitis A-Adult '
$wsAdult.Range('A3').Value = 'Hepatitis A-Adult '
$wsAdult.Range('A4').Value = 'Hepatitis A Adult '
$wsAdult.Range('A5').Value = 'Hepatitis A Adult '
$wsAdult.Range('A6').Value = 'Hepatitis A-Hepatitis B Adult '
$wsAdult.Range('A7').Value = 'Hepatitis A-Hepatitis B Adult '
$wsAdult.Range('A8').Value = 'Hepatitis B-Adult '
$wsAdult.Range('A9').Value = 'Hepatitis B-Adult '
$wsAdult.Range('A10').Value = 'Hepatitis B-Adult '
$wsAdult.Range('A11').Value = 'Hepatitis B-Adult '
$wsAdult.Range('A12').Value = 'HPV -Quadrivalent Human Papillomavirus Types 6, 11, 16 and 18 Recombinant Adult '
$wsAdult.Range('A13').Value = 'HPV-Human Papillomavirus 9 Valent '
$wsAdult.Range('A14').Value = 'HPV-Human Papillomavirus Bivalent Types 16 and 18 '
$wsAdult.Range('A15').Value = 'Measles, Mumps,  Rubella-Adult '
$wsAdult.Range('A16').Value = 'Meningococcal Conjugate '
$wsAdult.Range('A17').Value = 'Meningococcal Conjugate '
$wsAdult.Range('A21').Value = 'Pneumococcal 13-valent  (Adult)'
$wsAdult.Range('A24').Value = 'Tetanus and Diphtheria Toxoids'
$wsAdult.Range('A25').Value = 'Tetanus and Diphtheria Toxoids'
$wsAdult.Range('A30').Value = 'Varicella-Adult '
$wsPedFlu.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$wsAdultFlu.Range('A2').Value = 'Influenza  (Age 6 months and older)'
$wsPedFlu.Range('B2').Value = 'Fluzone Quadrivalent'
$wsAdultFlu.Range('B2').Value = 'Fluzone Quadrivalent'
$wsPedFlu.Range('A3').Value = 'Influenza  (Age 6-35 months)'
$wsPedFlu.Range('B3').Value = 'Fluzone Quadrivalent Pediatric dose No Preservative'
$wsPedFlu.Range('A4').Value = 'Influenza  (Age 36 months and older)'
$wsPedFlu.Range('A5').Value = 'Influenza  (Age 36 months and older)'
$wsPedFlu.Range('A6').Value = 'Influenza  (Age 36 months and older)'
$wsPedFlu.Range('A7').Value = 'Influenza  (Age 36 months and older)'
$wsAdultFlu.Range('A3').Value = 'Influenza  (Age 36 months and older)'
$wsAdultFlu.Range('A4').Value = 'Influenza  (Age 36 months and older)'
$wsAdultFlu.Range('A6').Value = 'Influenza  (Age 36 months and older)'
$wsAdultFlu.Range('A7').Value = 'Influenza  (Age 36 months and older)'
$wsPedFlu.Range('B4').Value = 'Fluzone Quadrivalent No-Preservative'
$wsPedFlu.Range('B5').Value = 'Fluzone Quadrivalent No-Preservative'
$wsPedFlu.Range('B6').Value = 'Fluarix Quadrivalent Preservative Free'
$wsAdultFlu.Range('B6').Value = 'Fluarix Quadrivalent Preservative Free'
$wsPedFlu.Range('B7').Value = 'FluLaval Quadrivalent'
$wsPedFlu.Range('A8').Value = 'Influenza  (Age 4 years and older)'
$wsAdultFlu.Range('A5').Value = 'Influenza  (Age 4 years and older)'
$wsPedFlu.Range('A9').Value = 'Influenza  Live, Intranasal (Age 2-49 years)'
$wsPedFlu.Range('B9').Value = 'FluMist Quadrivalent No Preservative'
$wsPedFlu.Range('A10').Value = 'Influenza  (Age 9 years and older)'
$wsPedFlu.Range('A11').Value = 'Influenza  (Age 9 years and older)'
$wsAdultFlu.Range('A8').Value = 'Influenza  (Age 9 years and older)'
$wsAdultFlu.Range('A9').Value = 'Influenza  (Age 9 years and older)'
$wsPedFlu.Range('B10').Value = 'Afluria No Preservative'
$wsAdultFlu.Range('B8').Value = 'Afluria No Preservative'
$wsPedFlu.Range('D10').Value = '10 pack-1 dose syringe'
$wsAdultFlu.Range('D8').Value = '10 pack-1 dose syringe'
$wsAdultFlu.Range('B3').Value = 'Fluzone Quadrivalent No Preservative'
$wsAdultFlu.Range('B4').Value = 'Fluzone Quadrivalent No Preservative'
$wsAdultFlu.Range('A10').Value = 'Influenza  (Age 18 years and older)'
